# Update 'F' column (想去人数) values across sheets per the source diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Cells.Item(4, 6).Value = 8367
$ws.Cells.Item(5, 6).Value = 8367
$ws.Cells.Item(8, 6).Value = 505
$ws.Cells.Item(9, 6).Value = 7336
$ws.Cells.Item(10, 6).Value = 588
$ws.Cells.Item(11, 6).Value = 510
$ws.Cells.Item(14, 6).Value = 165
$ws.Cells.Item(16, 6).Value = 169
$ws.Cells.Item(18, 6).Value = 144
$ws.Cells.Item(19, 6).Value = 12157
$ws.Cells.Item(20, 6).Value = 106
$ws.Cells.Item(22, 6).Value = 2462
$ws.Cells.Item(23, 6).Value = 3522
$ws.Cells.Item(26, 6).Value = 2906
$ws.Cells.Item(27, 6).Value = 109
$ws.Cells.Item(28, 6).Value = 112
$ws.Cells.Item(30, 6).Value = 3349
$ws.Cells.Item(32, 6).Value = 342
$ws.Cells.Item(33, 6).Value = 1713
$ws.Cells.Item(35, 6).Value = 132
$ws.Cells.Item(36, 6).Value = 6015
$ws.Cells.Item(38, 6).Value = 1830
$ws.Cells.Item(39, 6).Value = 1255
$ws.Cells.Item(40, 6).Value = 33
$ws.Cells.Item(41, 6).Value = 899
$ws.Cells.Item(43, 6).Value = 171
$ws.Cells.Item(45, 6).Value = 197
$ws.Cells.Item(46, 6).Value = 1124
$ws.Cells.Item(47, 6).Value = 1112
$ws.Cells.Item(48, 6).Value = 1581
$ws.Cells.Item(49, 6).Value = 19
$ws.Cells.Item(50, 6).Value = 116

$ws = $wb.Worksheets.Item("演出")
$ws.Cells.Item(10, 6).Value = 52
$ws.Cells.Item(22, 6).Value = 73

$ws = $wb.Worksheets.Item("本地生活")
$ws.Cells.Item(2, 6).Value = 320
$ws.Cells.Item(3, 6).Value = 465

$ws = $wb.Worksheets.Item("全部类型")
$ws.Cells.Item(5, 6).Value = 320
$ws.Cells.Item(8, 6).Value = 8367
$ws.Cells.Item(11, 6).Value = 506
$ws.Cells.Item(12, 6).Value = 7336
$ws.Cells.Item(13, 6).Value = 7336
$ws.Cells.Item(14, 6).Value = 588
$ws.Cells.Item(15, 6).Value = 510
$ws.Cells.Item(17, 6).Value = 165
$ws.Cells.Item(20, 6).Value = 169
$ws.Cells.Item(21, 6).Value = 144
$ws.Cells.Item(23, 6).Value = 12157
$ws.Cells.Item(24, 6).Value = 106
$ws.Cells.Item(27, 6).Value = 2462
$ws.Cells.Item(28, 6).Value = 2462
$ws.Cells.Item(29, 6).Value = 3522
$ws.Cells.Item(30, 6).Value = 109
$ws.Cells.Item(31, 6).Value = 112
$ws.Cells.Item(34, 6).Value = 3349
$ws.Cells.Item(36, 6).Value = 342
$ws.Cells.Item(37, 6).Value = 1713
$ws.Cells.Item(39, 6).Value = 132
$ws.Cells.Item(40, 6).Value = 6015
$ws.Cells.Item(41, 6).Value = 73
$ws.Cells.Item(42, 6).Value = 1830
$ws.Cells.Item(44, 6).Value = 1255
$ws.Cells.Item(45, 6).Value = 33
$ws.Cells.Item(46, 6).Value = 899
$ws.Cells.Item(47, 6).Value = 171
$ws.Cells.Item(48, 6).Value = 197
$ws.Cells.Item(49, 6).Value = 1124
$ws.Cells.Item(50, 6).Value = 1112
$ws.Cells.Item(51, 6).Value = 1581
$ws.Cells.Item(52, 6).Value = 116
